# "Navigate to Issue search, display all"
# Adds a new "browse_issues" worksheet after the existing "login_failed"
# sheet, populates it with the project-browse list, and leaves it as the
# active/selected sheet (matching the target workbook state).

$wb = $excel.ActiveWorkbook

# Existing sheet to anchor the new one after (the workbook's only/active
# sheet at the start of the edit).
$loginFailed = $wb.ActiveSheet

# Insert the new sheet right after "login_failed".
$browseIssues = $wb.Worksheets.Add($null, $loginFailed)
$browseIssues.Name = "browse_issues"

# Fill in the data.
$browseIssues.Range("A1").Value = "Projects to browse"
$browseIssues.Range("A2").Value = "TOUCAN"
$browseIssues.Range("A3").Value = "JETI"
$browseIssues.Range("A4").Value = "COALA"

# Leave selection/activation on the last entered cell of the new sheet,
# mirroring a user who just finished typing the list.
$browseIssues.Activate() | Out-Null
$browseIssues.Range("A4").Select() | Out-Null
